# v0.2.2 - Added Online/Presential
#
# - University: reduce the semester length ("Days") from 60 to 30.
# - Promotions: register the new "A3" promotion (already referenced by the
#   Subjects sheet) alongside A1/A2, using the same "TDA" grouping.
# - Rooms: drop the unused placeholder rooms L104-L106 and add a new
#   "Online" room (type "default") so subjects can be scheduled online
#   as well as in a physical ("Presential") room.

$wb = $excel.ActiveWorkbook

# --- Promotions: add the A3 / TDA column ---------------------------------
$wsPromotions = $wb.Worksheets.Item("Promotions")
$wsPromotions.Range("C1").Value = "A3"
$wsPromotions.Range("C2").Value = "TDA"
$wsPromotions.Activate()
[void]$wsPromotions.Range("D2").Select()

# --- Subjects: no data change, just record the last-used selection -------
$wsSubjects = $wb.Worksheets.Item("Subjects")
$wsSubjects.Activate()
[void]$wsSubjects.Range("A20").Select()

# --- Rooms: remove L104-L106, add the Online room -------------------------
$wsRooms = $wb.Worksheets.Item("Rooms")
$wsRooms.Range("A5:A7").EntireRow.Delete()
$wsRooms.Range("A5").Value = "Online"
$wsRooms.Range("B5").Value = "default"
$wsRooms.Activate()
[void]$wsRooms.Range("C5").Select()

# --- University: shorten the semester from 60 to 30 days ------------------
$wsUniversity = $wb.Worksheets.Item("University")
$wsUniversity.Range("B6").Value = 30
$wsUniversity.Activate()
